$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users")
$wsUserGroups = $wb.Worksheets.Item("User Groups")

# --- Users sheet: add new row for kyle.koyanagi ---
$wsUsers.Cells.Item(8, 1).Value = "kyle.koyanagi"
$wsUsers.Cells.Item(8, 3).Value = "Y"
$wsUsers.Cells.Item(8, 4).Formula = '=CONCATENATE("INSERT INTO AUTH_APP_USERS (",A$1, ", ", B$1, ", ", C$1, ") VALUES (''", A8, "'', ''", SUBSTITUTE(B8, "''", "''''"), "'', ''", C8, "'');")'

# --- User Groups sheet: add new row for kyle.koyanagi / DATA_ADMIN ---
$wsUserGroups.Cells.Item(8, 1).Value = "kyle.koyanagi"
$wsUserGroups.Cells.Item(8, 2).Value = "DATA_ADMIN"
$wsUserGroups.Cells.Item(8, 3).Formula = '=CONCATENATE("INSERT INTO AUTH_APP_USER_GROUPS (",A$1, ", ", B$1, ") VALUES ((SELECT APP_USER_ID FROM AUTH_APP_USERS WHERE APP_USER_NAME = ''", A8, "''), (SELECT APP_GROUP_ID FROM AUTH_APP_GROUPS WHERE APP_GROUP_CODE = ''", SUBSTITUTE(B8, "''", "''''"), "''));")'

# --- Update selections / active sheet to match the commit's final state ---
# User Groups tab loses focus; its selection moves to the newly added row.
$wsUserGroups.Range("A8").Select() | Out-Null

# Users tab becomes the active tab, with the new cell selected.
$wsUsers.Activate() | Out-Null
$wsUsers.Range("D8").Select() | Out-Null
